$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the student id in A2 (one digit shorter than before)
$ws.Range("A2").Value = 10152510288513

# Add a new row 3 with values 0..21 across columns A..V (as text)
$values = 0..21
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(3, $col).Value = "$($values[$i])"
}

# Update the active selection
$ws.Range("G9").Select()
